# "add nptg add overwrite option"
#
# The "Street" column (Landmark / Street pair) is removed from the Stops
# sheet's request template, the "Landmark" header is renamed to
# "Indicator", and the following columns all shift one place to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old "Street" column (column G) entirely - this shifts every
# column after it (StopType, GridType, Easting, Northing, Longitude,
# Latitude, LocalityCentre, CreationDateTime, ModificationDateTime,
# Modification, RevisionNumber, Status, AdministrativeAreaRef) one place
# to the left, and also shifts the associated data validation ranges.
$ws.Columns("G").Delete()

# The old "Landmark" column is now column F - rename its header to
# "Indicator" to reflect its new purpose.
$ws.Range("F1").Value = "Indicator"

# Reflect the new active selection (the cell under the renamed header).
$ws.Range("F2").Select() | Out-Null
